$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.141.82'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '1.859.07'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.79'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.81'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2844'
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06459'
$ws.Range("E10").Value = '  -2.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.94'
$ws.Range("E11").Value = '  -3.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07750'
$ws.Range("E12").Value = '  -3.42%  '
$ws.Range("D13").Value = '1.866.77'
$ws.Range("E13").Value = '  -0.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.33'
$ws.Range("E14").Value = '  -4.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6784'
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.035'
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '265.59'
$ws.Range("E17").Value = '  -1.99%  '
$ws.Range("D18").Value = '30.126.94'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.29'
$ws.Range("E19").Value = '  -5.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007571'
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '2.113.28'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.120'
$ws.Range("E24").Value = '  -3.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.085'
$ws.Range("E25").Value = '  -2.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.310'
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.05'
$ws.Range("E27").Value = '  -2.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.44'
$ws.Range("E28").Value = '  -2.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.877'
$ws.Range("E29").Value = '  -4.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.362'
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09903'
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.449'
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.199'
$ws.Range("E33").Value = '  -3.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.973'
$ws.Range("E34").Value = '  -2.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.04653'
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.112'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6853'
$ws.Range("E37").Value = '  -2.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.714'
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01828'
$ws.Range("E39").Value = '  -2.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.751'
$ws.Range("E40").Value = '  +3.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.275'
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.94'
$ws.Range("E42").Value = '  -2.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8303'
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.877'
$ws.Range("E45").Value = '  -4.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.93'
$ws.Range("E46").Value = '  -1.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4030'
$ws.Range("E47").Value = '  -3.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.085'
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '926.48'
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.917'
$ws.Range("E50").Value = '  -2.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.91'
